$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 36399
$ws.Cells.Item(2, 4).Value = 52682575
$ws.Cells.Item(3, 3).Value = 88352
$ws.Cells.Item(3, 4).Value = 129617155
$ws.Cells.Item(4, 3).Value = 30275
$ws.Cells.Item(4, 4).Value = 44864033
$ws.Cells.Item(5, 3).Value = 8393
$ws.Cells.Item(5, 4).Value = 12477676
$ws.Cells.Item(6, 3).Value = 1845
$ws.Cells.Item(6, 4).Value = 2743547
$ws.Cells.Item(7, 3).Value = 138
$ws.Cells.Item(7, 4).Value = 202093
$ws.Cells.Item(11, 3).Value = 39783
$ws.Cells.Item(11, 4).Value = 54060106
$ws.Cells.Item(12, 3).Value = 9343
$ws.Cells.Item(12, 4).Value = 13517559
$ws.Cells.Item(13, 3).Value = 25310
$ws.Cells.Item(13, 4).Value = 37138633
$ws.Cells.Item(14, 3).Value = 8107
$ws.Cells.Item(14, 4).Value = 12036168
$ws.Cells.Item(15, 3).Value = 2084
$ws.Cells.Item(15, 4).Value = 3101150
$ws.Cells.Item(16, 3).Value = 390
$ws.Cells.Item(16, 4).Value = 574123
$ws.Cells.Item(19, 3).Value = 9842
$ws.Cells.Item(19, 4).Value = 13070703
$ws.Cells.Item(20, 3).Value = 13049
$ws.Cells.Item(20, 4).Value = 18849366
$ws.Cells.Item(21, 3).Value = 30972
$ws.Cells.Item(21, 4).Value = 45481823
$ws.Cells.Item(22, 3).Value = 10030
$ws.Cells.Item(22, 4).Value = 14916231
$ws.Cells.Item(23, 3).Value = 2558
$ws.Cells.Item(23, 4).Value = 3807263
$ws.Cells.Item(24, 3).Value = 469
$ws.Cells.Item(24, 4).Value = 697592
$ws.Cells.Item(25, 3).Value = 34
$ws.Cells.Item(25, 4).Value = 50453
$ws.Cells.Item(26, 3).Value = 11353
$ws.Cells.Item(26, 4).Value = 15195673
$ws.Cells.Item(27, 3).Value = 7407
$ws.Cells.Item(27, 4).Value = 10731837
$ws.Cells.Item(28, 3).Value = 21974
$ws.Cells.Item(28, 4).Value = 32257013
$ws.Cells.Item(29, 3).Value = 7610
$ws.Cells.Item(29, 4).Value = 11325304
$ws.Cells.Item(30, 3).Value = 1910
$ws.Cells.Item(30, 4).Value = 2850544
$ws.Cells.Item(31, 3).Value = 341
$ws.Cells.Item(31, 4).Value = 508915
$ws.Cells.Item(33, 3).Value = 8072
$ws.Cells.Item(33, 4).Value = 10676850
$ws.Cells.Item(34, 3).Value = 3066
$ws.Cells.Item(34, 4).Value = 4421459
$ws.Cells.Item(35, 3).Value = 7460
$ws.Cells.Item(35, 4).Value = 10901189
$ws.Cells.Item(36, 3).Value = 3008
$ws.Cells.Item(36, 4).Value = 4456288
$ws.Cells.Item(37, 3).Value = 791
$ws.Cells.Item(37, 4).Value = 1179263
$ws.Cells.Item(38, 3).Value = 144
$ws.Cells.Item(38, 4).Value = 214232
$ws.Cells.Item(40, 3).Value = 2305
$ws.Cells.Item(40, 4).Value = 3112387
$ws.Cells.Item(41, 3).Value = 16718
$ws.Cells.Item(41, 4).Value = 24190852
$ws.Cells.Item(42, 3).Value = 49734
$ws.Cells.Item(42, 4).Value = 72943828
$ws.Cells.Item(43, 3).Value = 18548
$ws.Cells.Item(43, 4).Value = 27552530
$ws.Cells.Item(44, 3).Value = 5420
$ws.Cells.Item(44, 4).Value = 8072767
$ws.Cells.Item(45, 3).Value = 1104
$ws.Cells.Item(45, 4).Value = 1646792
$ws.Cells.Item(49, 3).Value = 16176
$ws.Cells.Item(49, 4).Value = 21572625
$ws.Cells.Item(50, 3).Value = 1864
$ws.Cells.Item(50, 4).Value = 2705440
$ws.Cells.Item(51, 3).Value = 6450
$ws.Cells.Item(51, 4).Value = 9489228
$ws.Cells.Item(52, 3).Value = 2232
$ws.Cells.Item(52, 4).Value = 3333324
$ws.Cells.Item(53, 3).Value = 726
$ws.Cells.Item(53, 4).Value = 1084305
$ws.Cells.Item(54, 3).Value = 170
$ws.Cells.Item(54, 4).Value = 251833
$ws.Cells.Item(56, 3).Value = 6102
$ws.Cells.Item(56, 4).Value = 8408761
$ws.Cells.Item(57, 3).Value = 810
$ws.Cells.Item(57, 4).Value = 1186584
$ws.Cells.Item(58, 3).Value = 2018
$ws.Cells.Item(58, 4).Value = 2994826
$ws.Cells.Item(59, 3).Value = 823
$ws.Cells.Item(59, 4).Value = 1226145
$ws.Cells.Item(60, 3).Value = 281
$ws.Cells.Item(60, 4).Value = 421258
$ws.Cells.Item(61, 3).Value = 69
$ws.Cells.Item(61, 4).Value = 103500
$ws.Cells.Item(63, 3).Value = 1190
$ws.Cells.Item(63, 4).Value = 1682894
$ws.Cells.Item(64, 3).Value = 14924
$ws.Cells.Item(64, 4).Value = 21566605
$ws.Cells.Item(65, 3).Value = 43679
$ws.Cells.Item(65, 4).Value = 63952099
$ws.Cells.Item(66, 3).Value = 15355
$ws.Cells.Item(66, 4).Value = 22828906
$ws.Cells.Item(67, 3).Value = 4445
$ws.Cells.Item(67, 4).Value = 6621560
$ws.Cells.Item(68, 3).Value = 886
$ws.Cells.Item(68, 4).Value = 1318596
$ws.Cells.Item(71, 3).Value = 14679
$ws.Cells.Item(71, 4).Value = 19387943
$ws.Cells.Item(72, 3).Value = 48892
$ws.Cells.Item(72, 4).Value = 71184838
$ws.Cells.Item(73, 3).Value = 140072
$ws.Cells.Item(73, 4).Value = 206454190
$ws.Cells.Item(74, 3).Value = 61193
$ws.Cells.Item(74, 4).Value = 91211890
$ws.Cells.Item(75, 3).Value = 19467
$ws.Cells.Item(75, 4).Value = 29090019
$ws.Cells.Item(76, 3).Value = 4459
$ws.Cells.Item(76, 4).Value = 6662472
$ws.Cells.Item(83, 3).Value = 48504
$ws.Cells.Item(83, 4).Value = 66201533
$ws.Cells.Item(84, 3).Value = 4389
$ws.Cells.Item(84, 4).Value = 6362228
$ws.Cells.Item(85, 3).Value = 11120
$ws.Cells.Item(85, 4).Value = 16342446
$ws.Cells.Item(86, 3).Value = 3772
$ws.Cells.Item(86, 4).Value = 5621415
$ws.Cells.Item(87, 3).Value = 1312
$ws.Cells.Item(87, 4).Value = 1960489
$ws.Cells.Item(91, 3).Value = 5107
$ws.Cells.Item(91, 4).Value = 6878553
$ws.Cells.Item(92, 3).Value = 1487
$ws.Cells.Item(92, 4).Value = 2149412
$ws.Cells.Item(93, 3).Value = 4885
$ws.Cells.Item(93, 4).Value = 7195584
$ws.Cells.Item(94, 3).Value = 1864
$ws.Cells.Item(94, 4).Value = 2777499
$ws.Cells.Item(95, 3).Value = 659
$ws.Cells.Item(95, 4).Value = 987460
$ws.Cells.Item(99, 3).Value = 3302
$ws.Cells.Item(99, 4).Value = 4383021
$ws.Cells.Item(100, 3).Value = 565
$ws.Cells.Item(100, 4).Value = 842964
$ws.Cells.Item(101, 3).Value = 321
$ws.Cells.Item(101, 4).Value = 479630
$ws.Cells.Item(102, 3).Value = 113
$ws.Cells.Item(102, 4).Value = 169500
$ws.Cells.Item(105, 3).Value = 10451
$ws.Cells.Item(105, 4).Value = 15182390
$ws.Cells.Item(106, 3).Value = 28575
$ws.Cells.Item(106, 4).Value = 41993886
$ws.Cells.Item(107, 3).Value = 9573
$ws.Cells.Item(107, 4).Value = 14236157
$ws.Cells.Item(108, 3).Value = 2612
$ws.Cells.Item(108, 4).Value = 3894307
$ws.Cells.Item(109, 3).Value = 462
$ws.Cells.Item(109, 4).Value = 690482
$ws.Cells.Item(112, 3).Value = 9513
$ws.Cells.Item(112, 4).Value = 12592630
$ws.Cells.Item(113, 3).Value = 29479
$ws.Cells.Item(113, 4).Value = 42543771
$ws.Cells.Item(114, 3).Value = 64538
$ws.Cells.Item(114, 4).Value = 94491545
$ws.Cells.Item(115, 3).Value = 20884
$ws.Cells.Item(115, 4).Value = 31047979
$ws.Cells.Item(116, 3).Value = 5872
$ws.Cells.Item(116, 4).Value = 8749280
$ws.Cells.Item(117, 3).Value = 1067
$ws.Cells.Item(117, 4).Value = 1595006
$ws.Cells.Item(118, 3).Value = 68
$ws.Cells.Item(118, 4).Value = 99420
$ws.Cells.Item(121, 3).Value = 24990
$ws.Cells.Item(121, 4).Value = 33422976
$ws.Cells.Item(122, 3).Value = 34569
$ws.Cells.Item(122, 4).Value = 49930096
$ws.Cells.Item(123, 3).Value = 74367
$ws.Cells.Item(123, 4).Value = 108823337
$ws.Cells.Item(124, 3).Value = 23146
$ws.Cells.Item(124, 4).Value = 34364029
$ws.Cells.Item(125, 3).Value = 6166
$ws.Cells.Item(125, 4).Value = 9164367
$ws.Cells.Item(126, 3).Value = 1151
$ws.Cells.Item(126, 4).Value = 1710911
$ws.Cells.Item(130, 3).Value = 30480
$ws.Cells.Item(130, 4).Value = 40548885
$ws.Cells.Item(131, 3).Value = 12866
$ws.Cells.Item(131, 4).Value = 18631818
$ws.Cells.Item(132, 3).Value = 31590
$ws.Cells.Item(132, 4).Value = 46421951
$ws.Cells.Item(133, 3).Value = 11243
$ws.Cells.Item(133, 4).Value = 16706059
$ws.Cells.Item(134, 3).Value = 2880
$ws.Cells.Item(134, 4).Value = 4294581
$ws.Cells.Item(135, 3).Value = 463
$ws.Cells.Item(135, 4).Value = 688490
$ws.Cells.Item(138, 3).Value = 10515
$ws.Cells.Item(138, 4).Value = 14044739
$ws.Cells.Item(139, 3).Value = 33807
$ws.Cells.Item(139, 4).Value = 48851883
$ws.Cells.Item(140, 3).Value = 78820
$ws.Cells.Item(140, 4).Value = 115526341
$ws.Cells.Item(141, 3).Value = 23685
$ws.Cells.Item(141, 4).Value = 35211145
$ws.Cells.Item(142, 3).Value = 6177
$ws.Cells.Item(142, 4).Value = 9220086
$ws.Cells.Item(143, 3).Value = 1362
$ws.Cells.Item(143, 4).Value = 2026192
$ws.Cells.Item(146, 3).Value = 28222
$ws.Cells.Item(146, 4).Value = 38171616
